# Applies the "daily scrum, otletek dokumentalasa" edit:
#  - title paragraph gets a trailing ":" run, both sz=48/szCs=48
#  - new "Product owner utasitasok:" heading paragraph (sz=32/szCs=32)
#  - existing body paragraphs (Osztalykep .. -leiras a tantargyakrol) are untouched
#  - the trailing " " paragraph loses the _GoBack bookmark
#  - new "Fooldal otletek:" heading + 5 bullet paragraphs are appended, with the
#    _GoBack bookmark now wrapping the end of the last ("-footer (talan)") paragraph
#  - a final empty paragraph is appended after it
#
# Implemented via a single Range.InsertXML() call over the whole document body so
# every run/proofErr/bookmark/pPr element can be reproduced exactly as authored.

$d = $word.ActiveDocument

$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>11.E főoldal</w:t></w:r><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>utasitasok</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Osztálykép</w:t></w:r></w:p><w:p><w:r><w:t>Kis leírás</w:t></w:r></w:p><w:p><w:r><w:t>Órarend:</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>-kattintható tantárgyak</w:t></w:r></w:p><w:p><w:r><w:t>Tantárgy weblapok, benne:</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>-kötelező animáció</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>-tanárok</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>-leírás a tantárgyakról</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Fooldal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>otletek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>navbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fooldalon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>belul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> vezet helyekre)</w:t></w:r></w:p><w:p><w:r><w:t>-banner (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>osztalykep</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">-kis </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>leiras</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>orarend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> az </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>aljan</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>footer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (talan)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full = $d.Content
[void]$full.InsertXML($xmlPayload)

Write-Output ("paragraphs=" + $d.Paragraphs.Count)
